$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values get swapped between the paired rows.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Row pairs to swap.
$pairs = @(
    @(3, 4),
    @(18, 19),
    @(20, 21)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
